$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'24.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.212"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05785"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'6.519"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'3.124"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8159"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.8524"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1359"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.06976"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03148"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.02878"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09383"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.745"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001527"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04701"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.0005985"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.006288"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.001235"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.004534"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'0.00008601"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'3.498"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Value = "'0.3175"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'0.1339"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'0.1327"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "'0.0002331"
$ws.Range("D28").Style = "Normal"
$ws.Range("D40").Value = "'0.03647"
$ws.Range("D40").Style = "Normal"
$ws.Range("B41").Value = "'KickToken"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'0.006301"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'40KickTokenKICK"
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'BKEXToken"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.1054"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'41BKEXTokenBKK"
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'CEJI"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'0.002903"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'42CEJICEJIBestin24h"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.007472"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005274"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'0.3113"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.002336"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'47BOLOBOLOWorstin24h"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002002"
$ws.Range("D50").Style = "Normal"
